$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current last row (row 11), pushing it down to row 13.
$ws.Rows.Item(11).Insert(-4121)  # xlShiftDown
$ws.Rows.Item(11).Insert(-4121)  # xlShiftDown

# --- Column A (Identificador) for the two new rows: copy format from A10 ---
$ws.Range("A10").Copy()
$ws.Range("A11:A12").PasteSpecial(-4122)  # xlPasteFormats

# --- Column B (Descrição Resumida) for the two new rows: copy format from B10 ---
$ws.Range("B10").Copy()
$ws.Range("B11:B12").PasteSpecial(-4122)  # xlPasteFormats

# Remove the bottom border on B11/B12 since they are no longer the last row
$ws.Range("B11").Borders.Item(9).LineStyle = -4142  # xlEdgeBottom -> xlLineStyleNone
$ws.Range("B12").Borders.Item(9).LineStyle = -4142  # xlEdgeBottom -> xlLineStyleNone

$excel.CutCopyMode = 0

# --- Set the cell values (order chosen to match shared-string append order) ---
# Row 12 (new): RNF-11 / Permitir acesso ao sistema através de sessão única
$ws.Range("A12").Value = "RNF-11"
$ws.Range("B12").Value = "Permitir acesso ao sistema através de sessão única"

# Row 11: RNF-10 keeps its identifier, but description becomes "Validar entradas de usuários"
$ws.Range("A11").Value = "RNF-10"
$ws.Range("B11").Value = "Validar entradas de usuários"

# Row 13 (previously row 11, pushed down): RNF-12 / original "Armazenar..." text
$ws.Range("A13").Value = "RNF-12"
$ws.Range("B13").Value = "Armazenar os dados por período indeterminado de tempo"

# Restore the dimension / selection bookkeeping to match the authored workbook
$ws.Range("B31").Select()
